$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value2 = 2684.0417
$ws.Range("J17").Value2 = 2719.8723
$ws.Range("L17").Value2 = 8159.6169
$ws.Range("N17").Value2 = -8495.616900000001

$ws.Range("H38").Value2 = 433
$ws.Range("I38").Value2 = 433
$ws.Range("J38").Value2 = 0
$ws.Range("K38").Value2 = 1299
$ws.Range("L38").Value2 = 0
$ws.Range("M38").Value2 = -927
$ws.Range("N38").ClearContents()

$ws.Range("H137").Value2 = 1645.1666
$ws.Range("I137").Value2 = 1215.9286
$ws.Range("J137").Value2 = 2646.7222
$ws.Range("K137").Value2 = 3647.7858
$ws.Range("L137").Value2 = 7940.1666
$ws.Range("M137").Value2 = -1097.7858
$ws.Range("N137").Value2 = -13040.1666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value2 = 1348.4108
$ws.Range("I45").Value2 = 1288.2291
$ws.Range("J45").Value2 = 1709.5
$ws.Range("K45").Value2 = 1288.2291
$ws.Range("L45").Value2 = 1709.5
$ws.Range("M45").Value2 = -911.2291
$ws.Range("N45").Value2 = -2463.5

$ws.Range("H61").Value2 = 7532.8936
$ws.Range("I61").Value2 = 3982.6052
$ws.Range("J61").Value2 = 22523
$ws.Range("K61").Value2 = 3982.6052
$ws.Range("L61").Value2 = 22523
$ws.Range("M61").Value2 = -3770.6052
$ws.Range("N61").Value2 = -22947

$ws.Range("H74").Value2 = 5350
$ws.Range("I74").Value2 = 2349.318
$ws.Range("J74").Value2 = 13601.875
$ws.Range("K74").Value2 = 2349.318
$ws.Range("L74").Value2 = 13601.875
$ws.Range("M74").Value2 = -1475.318
$ws.Range("N74").Value2 = -15349.875

$ws.Range("H77").Value2 = 5350
$ws.Range("I77").Value2 = 2349.318
$ws.Range("J77").Value2 = 13601.875
$ws.Range("K77").Value2 = 11746.59
$ws.Range("L77").Value2 = 68009.375
$ws.Range("M77").Value2 = -7378.59
$ws.Range("N77").Value2 = -76745.375

$ws.Range("H136").Value2 = 7532.8936
$ws.Range("I136").Value2 = 3982.6052
$ws.Range("J136").Value2 = 22523
$ws.Range("K136").Value2 = 11947.8156
$ws.Range("L136").Value2 = 67569
$ws.Range("M136").Value2 = -9397.8156
$ws.Range("N136").Value2 = -72669

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H14").Value2 = 10000
$ws.Range("J14").Value2 = 10000
$ws.Range("L14").Value2 = 10000
$ws.Range("N14").Value2 = -10344

$ws.Range("H35").Value2 = 36073
$ws.Range("J35").Value2 = 36073
$ws.Range("L35").Value2 = 36073
$ws.Range("N35").Value2 = -36693

$ws.Range("H86").Value2 = 1790.4517
$ws.Range("I86").Value2 = 1717.138
$ws.Range("J86").Value2 = 2853.5
$ws.Range("K86").Value2 = 1717.138
$ws.Range("L86").Value2 = 2853.5
$ws.Range("M86").Value2 = -594.1379999999999
$ws.Range("N86").Value2 = -5099.5

$ws.Range("H88").Value2 = 39900
$ws.Range("J88").Value2 = 39900
$ws.Range("L88").Value2 = 39900
$ws.Range("N88").Value2 = -40712

$ws.Range("H89").Value2 = 1790.4517
$ws.Range("I89").Value2 = 1717.138
$ws.Range("J89").Value2 = 2853.5
$ws.Range("K89").Value2 = 8585.689999999999
$ws.Range("L89").Value2 = 14267.5
$ws.Range("M89").Value2 = -2969.689999999999
$ws.Range("N89").Value2 = -25499.5

$ws.Range("H91").Value2 = 39900
$ws.Range("J91").Value2 = 39900
$ws.Range("L91").Value2 = 39900
$ws.Range("N91").Value2 = -42708

$ws.Range("H95").Value2 = 19966.666
$ws.Range("J95").Value2 = 19966.666
$ws.Range("L95").Value2 = 19966.666
$ws.Range("N95").Value2 = -25458.666

$ws.Range("H134").Value2 = 30082.584
$ws.Range("I134").Value2 = 2093.24
$ws.Range("K134").Value2 = 6279.719999999999
$ws.Range("M134").Value2 = -3744.719999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 1889.1167
$ws.Range("I31").Value2 = 1296.14
$ws.Range("J31").Value2 = 4854
$ws.Range("K31").Value2 = 1296.14
$ws.Range("L31").Value2 = 4854
$ws.Range("M31").Value2 = -1001.14
$ws.Range("N31").Value2 = -5444

$ws.Range("H34").Value2 = 1889.1167
$ws.Range("I34").Value2 = 1296.14
$ws.Range("J34").Value2 = 4854
$ws.Range("K34").Value2 = 1296.14
$ws.Range("L34").Value2 = 4854
$ws.Range("M34").Value2 = -1094.14
$ws.Range("N34").Value2 = -5258

$ws.Range("H59").Value2 = 29500
$ws.Range("J59").Value2 = 29500
$ws.Range("L59").Value2 = 29500
$ws.Range("N59").Value2 = -31790

$ws.Range("H88").Value2 = 30245.545
$ws.Range("J88").Value2 = 37837.625
$ws.Range("L88").Value2 = 37837.625
$ws.Range("N88").Value2 = -38649.625

$ws.Range("H91").Value2 = 30245.545
$ws.Range("J91").Value2 = 37837.625
$ws.Range("L91").Value2 = 37837.625
$ws.Range("N91").Value2 = -40645.625

$ws.Range("H132").Value2 = 3210.0435
$ws.Range("I132").Value2 = 3490.8367
$ws.Range("J132").Value2 = 2522.1
$ws.Range("K132").Value2 = 10472.5101
$ws.Range("L132").Value2 = 7566.299999999999
$ws.Range("M132").Value2 = -7942.5101
$ws.Range("N132").Value2 = -12626.3

$ws.Range("H134").Value2 = 2789.2856
$ws.Range("I134").Value2 = 1772.2084
$ws.Range("K134").Value2 = 5316.6252
$ws.Range("M134").Value2 = -2781.6252

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value2 = 1177
$ws.Range("J107").Value2 = 1696.5
$ws.Range("L107").Value2 = 5089.5
$ws.Range("N107").Value2 = -8929.5

$ws.Range("H131").Value2 = 1391.1063
$ws.Range("J131").Value2 = 1167.6666
$ws.Range("L131").Value2 = 3502.9998
$ws.Range("N131").Value2 = -13582.9998

$ws.Range("H132").Value2 = 1604.28
$ws.Range("I132").Value2 = 1671.5
$ws.Range("J132").Value2 = 1542.2307
$ws.Range("K132").Value2 = 15043.5
$ws.Range("L132").Value2 = 13880.0763
$ws.Range("M132").Value2 = -12513.5
$ws.Range("N132").Value2 = -18940.0763

$ws.Range("H140").Value2 = 2324.0476
$ws.Range("I140").Value2 = 1893.0625
$ws.Range("J140").Value2 = 3703.2
$ws.Range("K140").Value2 = 5679.1875
$ws.Range("L140").Value2 = 11109.6
$ws.Range("M140").Value2 = -499.1875
$ws.Range("N140").Value2 = -21469.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value2 = 34284.285
$ws.Range("J39").Value2 = 34284.285
$ws.Range("L39").Value2 = 34284.285
$ws.Range("N39").Value2 = -35348.285

$ws.Range("H97").Value2 = 1522.2
$ws.Range("I97").Value2 = 2202.5
$ws.Range("K97").Value2 = 2202.5
$ws.Range("M97").Value2 = -1706.5

$ws.Range("H132").Value2 = 6280.9165
$ws.Range("I132").Value2 = 2268.762
$ws.Range("K132").Value2 = 6806.286
$ws.Range("M132").Value2 = -4276.286

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value2 = 45000
$ws.Range("J6").Value2 = 45000
$ws.Range("L6").Value2 = 45000
$ws.Range("N6").Value2 = -45224

$ws.Range("H16").Value2 = 1345.591
$ws.Range("I16").Value2 = 928.06665
$ws.Range("J16").Value2 = 2240.2856
$ws.Range("K16").Value2 = 928.06665
$ws.Range("L16").Value2 = 2240.2856
$ws.Range("M16").Value2 = -758.06665
$ws.Range("N16").Value2 = -2580.2856

$ws.Range("H95").Value2 = 30172
$ws.Range("J95").Value2 = 30172
$ws.Range("L95").Value2 = 30172
$ws.Range("N95").Value2 = -35664

$ws.Range("H132").Value2 = 2352.0962
$ws.Range("I132").Value2 = 2137.587
$ws.Range("J132").Value2 = 3996.6667
$ws.Range("K132").Value2 = 6412.761
$ws.Range("L132").Value2 = 11990.0001
$ws.Range("M132").Value2 = -3882.761
$ws.Range("N132").Value2 = -17050.0001

$ws.Range("H136").Value2 = 3450.2769
$ws.Range("I136").Value2 = 1780.9149
$ws.Range("J136").Value2 = 7809.1665
$ws.Range("K136").Value2 = 5342.7447
$ws.Range("L136").Value2 = 23427.4995
$ws.Range("M136").Value2 = -2792.7447
$ws.Range("N136").Value2 = -28527.4995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value2 = 33725
$ws.Range("J92").Value2 = 32966.668
$ws.Range("L92").Value2 = 32966.668
$ws.Range("N92").Value2 = -37958.668

$ws.Range("H94").Value2 = 22666.666
$ws.Range("J94").Value2 = 22666.666
$ws.Range("L94").Value2 = 22666.666
$ws.Range("N94").Value2 = -24468.666

$ws.Range("H95").Value2 = 40000
$ws.Range("J95").Value2 = 40000
$ws.Range("L95").Value2 = 40000
$ws.Range("N95").Value2 = -45492

$ws.Range("H136").Value2 = 4114.617
$ws.Range("I136").Value2 = 2635.1724
$ws.Range("J136").Value2 = 7845.391
$ws.Range("K136").Value2 = 7905.5172
$ws.Range("L136").Value2 = 23536.173
$ws.Range("M136").Value2 = -5355.5172
$ws.Range("N136").Value2 = -28636.173
